$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.393.30"
$ws.Cells.Item(2, 5).Value = "  +0.09%  "

$ws.Cells.Item(3, 4).Value = "1.846.70"
$ws.Cells.Item(3, 5).Value = "  +0.05%  "

$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "239.83"
$ws.Cells.Item(5, 5).Value = "  -0.14%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.6293"
$ws.Cells.Item(6, 5).Value = "  -0.99%  "

$ws.Cells.Item(7, 5).Value = "  +0.01%  "

$dCell = $ws.Cells.Item(8, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.07593"
$ws.Cells.Item(8, 5).Value = "  +0.59%  "

$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.2927"
$ws.Cells.Item(9, 5).Value = "  -1.26%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "24.46"
$ws.Cells.Item(10, 5).Value = "  -0.64%  "

$ws.Cells.Item(11, 5).Value = "  +0.04%  "

$ws.Cells.Item(12, 4).Value = "1.843.97"
$ws.Cells.Item(12, 5).Value = "  -7.11%  "

$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001095"
$ws.Cells.Item(13, 5).Value = "  +10.21%  "

$ws.Cells.Item(14, 5).Value = "  +0.24%  "

$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.6780"
$ws.Cells.Item(15, 5).Value = "  -0.77%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "83.65"
$ws.Cells.Item(16, 5).Value = "  +0.70%  "

$ws.Cells.Item(17, 4).Value = "2.095.55"
$ws.Cells.Item(17, 5).Value = "  -7.45%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.159"
$ws.Cells.Item(18, 5).Value = "  -0.08%  "

$ws.Cells.Item(19, 4).Value = "29.413.97"

$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "228.66"
$ws.Cells.Item(20, 5).Value = "  -0.52%  "

$ws.Cells.Item(22, 5).Value = "  +0.04%  "

$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "7.423"
$ws.Cells.Item(23, 5).Value = "  -1.79%  "

$ws.Cells.Item(24, 5).Value = "  +0.03%  "

$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "157.42"
$ws.Cells.Item(25, 5).Value = "  +0.26%  "

$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.1395"
$ws.Cells.Item(26, 5).Value = "  -0.53%  "

$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "8.361"
$ws.Cells.Item(27, 5).Value = "  -0.10%  "

$dCell = $ws.Cells.Item(28, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "17.61"
$ws.Cells.Item(28, 5).Value = "  -0.19%  "

$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.463"
$ws.Cells.Item(29, 5).Value = "  -0.11%  "

$dCell = $ws.Cells.Item(30, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.298"
$ws.Cells.Item(30, 5).Value = "  +4.09%  "

$dCell = $ws.Cells.Item(31, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.05604"
$ws.Cells.Item(31, 5).Value = "  -1.78%  "

$ws.Cells.Item(32, 5).Value = "  -0.55%  "

$ws.Cells.Item(33, 5).Value = "  +0.14%  "

$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.848"
$ws.Cells.Item(34, 5).Value = "  +0.20%  "

$ws.Cells.Item(35, 5).Value = "  +0.05%  "

$ws.Cells.Item(36, 5).Value = "  -1.05%  "

$ws.Cells.Item(37, 5).Value = "  -0.24%  "

$ws.Cells.Item(38, 4).Value = "1.229.69"
$ws.Cells.Item(38, 5).Value = "  -1.98%  "

$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.01797"
$ws.Cells.Item(39, 5).Value = "  -0.44%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.767"
$ws.Cells.Item(40, 5).Value = "  -0.85%  "

$dCell = $ws.Cells.Item(41, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.444"
$ws.Cells.Item(41, 5).Value = "  +4.14%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.9069"
$ws.Cells.Item(42, 5).Value = "  -0.07%  "

$ws.Cells.Item(43, 5).Value = "  -0.06%  "

$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "101.79"
$ws.Cells.Item(44, 5).Value = "  +0.04%  "

$dCell = $ws.Cells.Item(45, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "66.04"
$ws.Cells.Item(45, 5).Value = "  -0.62%  "

$ws.Cells.Item(46, 5).Value = "  +3.45%  "

$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "7.189"
$ws.Cells.Item(47, 5).Value = "  +1.98%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.4019"
$ws.Cells.Item(48, 5).Value = "  -0.14%  "

$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "8.975"
$ws.Cells.Item(49, 5).Value = "  -2.34%  "

$ws.Cells.Item(50, 5).Value = "  -1.59%  "

$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.1121"
$ws.Cells.Item(51, 5).Value = "  -0.62%  "
